$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H28").Value = 376.70834
$ws.Range("J28").Value = 441
$ws.Range("L28").Value = 441
$ws.Range("N28").Value = -1411

$ws.Range("H33").Value = 218.77777
$ws.Range("I33").Value = 218.77777
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 218.77777
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = 10.22223

$ws.Range("H112").Value = 1853.7368
$ws.Range("J112").Value = 1983.8125
$ws.Range("L112").Value = 5951.4375
$ws.Range("N112").Value = -8167.4375

$ws.Range("H132").Value = 1807.8644
$ws.Range("I132").Value = 1855.5962
$ws.Range("J132").Value = 1453.2858
$ws.Range("K132").Value = 5566.7886
$ws.Range("L132").Value = 4359.857400000001
$ws.Range("M132").Value = -3036.7886
$ws.Range("N132").Value = -9419.857400000001

$ws.Range("H137").Value = 2081
$ws.Range("I137").Value = 1977.5
$ws.Range("J137").Value = 2115.5
$ws.Range("K137").Value = 5932.5
$ws.Range("L137").Value = 6346.5
$ws.Range("M137").Value = -3382.5
$ws.Range("N137").Value = -11446.5

$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 6625.091
$ws.Range("I32").Value = 6093.9287
$ws.Range("J32").Value = 9599.6
$ws.Range("K32").Value = 6093.9287
$ws.Range("L32").Value = 9599.6
$ws.Range("M32").Value = -5806.9287
$ws.Range("N32").Value = -10173.6

$ws.Range("H37").Value = 4189128.2
$ws.Range("J37").Value = 24989.445
$ws.Range("L37").Value = 24989.445
$ws.Range("N37").Value = -25535.445

$ws.Range("H45").Value = 7091.5
$ws.Range("I45").Value = 14688.375
$ws.Range("K45").Value = 14688.375
$ws.Range("M45").Value = -14311.375

$ws.Range("H110").Value = 1376.25
$ws.Range("I110").Value = 1434
$ws.Range("K110").Value = 1434
$ws.Range("M110").Value = 611

$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("N113").Value = 0

$ws.Range("H122").Value = 1986.2122
$ws.Range("I122").Value = 1357.3636
$ws.Range("J122").Value = 3243.9092
$ws.Range("K122").Value = 4072.0908
$ws.Range("L122").Value = 9731.7276
$ws.Range("M122").Value = -1622.0908
$ws.Range("N122").Value = -14631.7276

$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 2084.8
$ws.Range("J16").Value = 1799
$ws.Range("L16").Value = 1799
$ws.Range("N16").Value = -2373

$ws.Range("H31").Value = 11458
$ws.Range("I31").Value = 1793.3077
$ws.Range("K31").Value = 1793.3077
$ws.Range("M31").Value = -1498.3077

$ws.Range("H34").Value = 11458
$ws.Range("I34").Value = 1793.3077
$ws.Range("K34").Value = 1793.3077
$ws.Range("M34").Value = -1591.3077

$ws.Range("H58").Value = 2203.6667
$ws.Range("I58").Value = 1408
$ws.Range("J58").Value = 2999.3333
$ws.Range("K58").Value = 1408
$ws.Range("L58").Value = 2999.3333
$ws.Range("M58").Value = -1205
$ws.Range("N58").Value = -3405.3333

$ws.Range("H113").Value = 2084.8
$ws.Range("J113").Value = 1799
$ws.Range("L113").Value = 1799
$ws.Range("N113").Value = -6139

$ws.Range("H136").Value = 2203.6667
$ws.Range("I136").Value = 1408
$ws.Range("J136").Value = 2999.3333
$ws.Range("K136").Value = 4224
$ws.Range("L136").Value = 8997.999899999999
$ws.Range("M136").Value = -1674
$ws.Range("N136").Value = -14097.9999

$ws = $wb.Worksheets.Item(5)
$ws.Range("H6").Value = 97
$ws.Range("I6").Value = 97
$ws.Range("K6").Value = 291
$ws.Range("M6").Value = -178

$ws.Range("H23").Value = 111.92308
$ws.Range("I23").Value = 146.66667
$ws.Range("K23").Value = 440.00001
$ws.Range("M23").Value = -205.00001

$ws.Range("H44").Value = 102.5
$ws.Range("I44").Value = 102.5
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 307.5
$ws.Range("L44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = 90.5

$ws.Range("H107").Value = 658.975
$ws.Range("J107").Value = 656.32434
$ws.Range("L107").Value = 1968.97302
$ws.Range("N107").Value = -5808.973019999999

$ws.Range("H123").Value = 880
$ws.Range("I123").Value = 880
$ws.Range("K123").Value = 2640
$ws.Range("M123").Value = -190

$ws.Range("H131").Value = 6820.0527
$ws.Range("I131").Value = 12514.333
$ws.Range("J131").Value = 1695.2
$ws.Range("K131").Value = 37542.999
$ws.Range("L131").Value = 5085.6
$ws.Range("M131").Value = -32502.999
$ws.Range("N131").Value = -15165.6

$ws = $wb.Worksheets.Item(6)
$ws.Range("H102").Value = 2711.8333
$ws.Range("I102").Value = 2822.682
$ws.Range("K102").Value = 2822.682
$ws.Range("M102").Value = -1200.682

$ws.Range("H113").Value = 3982.44
$ws.Range("I113").Value = 3507.4614
$ws.Range("J113").Value = 4497
$ws.Range("K113").Value = 3507.4614
$ws.Range("L113").Value = 4497
$ws.Range("M113").Value = -1337.4614
$ws.Range("N113").Value = -8837

$ws.Range("H126").Value = 1810.9
$ws.Range("I126").Value = 1489
$ws.Range("J126").Value = 3098.5
$ws.Range("K126").Value = 4467
$ws.Range("L126").Value = 9295.5
$ws.Range("M126").Value = -1997
$ws.Range("N126").Value = -14235.5

$ws = $wb.Worksheets.Item(7)
$ws.Range("H61").Value = 36555.95
$ws.Range("I61").Value = 1641
$ws.Range("K61").Value = 1641
$ws.Range("M61").Value = -1439

$ws.Range("H63").Value = 23085
$ws.Range("J63").Value = 23085
$ws.Range("L63").Value = 23085
$ws.Range("N63").Value = -24583

$ws.Range("H66").Value = 23085
$ws.Range("J66").Value = 23085
$ws.Range("L66").Value = 69255
$ws.Range("N66").Value = -76743

$ws.Range("H76").Value = 19733
$ws.Range("J76").Value = 19733
$ws.Range("L76").Value = 19733
$ws.Range("N76").Value = -20409

$ws.Range("H79").Value = 19733
$ws.Range("J79").Value = 19733
$ws.Range("L79").Value = 19733
$ws.Range("N79").Value = -22073

$ws.Range("H103").Value = 27798.5
$ws.Range("J103").Value = 27798.5
$ws.Range("L103").Value = 27798.5
$ws.Range("N103").Value = -30142.5

$ws.Range("H113").Value = 36555.95
$ws.Range("I113").Value = 1641
$ws.Range("K113").Value = 1641
$ws.Range("M113").Value = 529

$ws = $wb.Worksheets.Item(8)
$ws.Range("H82").Value = 28860
$ws.Range("J82").Value = 28860
$ws.Range("L82").Value = 28860
$ws.Range("N82").Value = -29626

$ws.Range("H85").Value = 28860
$ws.Range("J85").Value = 28860
$ws.Range("L85").Value = 28860
$ws.Range("N85").Value = -31512

$ws.Range("H101").Value = 19994.5
$ws.Range("J101").Value = 19994.5
$ws.Range("L101").Value = 19994.5
$ws.Range("N101").Value = -26484.5

$ws.Range("H112").Value = 39057
$ws.Range("J112").Value = 39057
$ws.Range("L112").Value = 39057
$ws.Range("N112").Value = -42011

$ws.Range("H113").Value = 223.26666
$ws.Range("I113").Value = 198.04546
$ws.Range("K113").Value = 594.1363799999999
$ws.Range("M113").Value = 1575.86362
